# Update column F (dSF) values to repull data / push all data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = -4
$ws.Range("F10").Value = 3
$ws.Range("F14").Value = -1
$ws.Range("F15").Value = -13
$ws.Range("F17").Value = -8
$ws.Range("F19").Value = -5
$ws.Range("F20").Value = -7
$ws.Range("F23").Value = -5
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 4
